$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing columns D:K to E:L.
# This makes room for a new (most recent) reporting period in the
# Income Statement / Balance Sheet / Cash Flow Statement blocks.
$ws.Columns("D:D").Insert(-4161)

# Carry the number formatting (date format for header rows, #,##0 for
# the data rows) from the column that used to be D (now E) into the
# freshly inserted column D.
$ws.Columns("E:E").Copy()
$ws.Columns("D:D").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New period data for column D (row number -> value; $null -> leave blank).
$colD = @{
    7 = 43373
    8 = 3743000
    9 = 2656000
    10 = 1087000
    11 = $null
    12 = 85000
    13 = 0
    14 = 74000
    15 = 95000
    16 = $null
    17 = 3508000
    18 = 235000
    19 = $null
    20 = 15000
    21 = 562000
    22 = 136000
    23 = 114000
    24 = 9000
    25 = 0
    26 = 105000
    27 = 105000
    28 = 0
    29 = 9000
    30 = 0
    31 = 0
    32 = -15000
    33 = 114000
    34 = 0
    35 = 114000
    38 = 43373
    39 = $null
    40 = $null
    41 = 294000
    42 = 0
    43 = 681000
    44 = 663000
    45 = 74000
    46 = 1712000
    47 = 460000
    48 = 1899000
    49 = 3689000
    50 = 0
    51 = 0
    52 = 492000
    53 = 0
    54 = 8252000
    55 = $null
    56 = $null
    57 = 483000
    58 = 254000
    59 = 338000
    60 = 1075000
    61 = 2275000
    62 = 1496000
    63 = 0
    64 = 0
    65 = 0
    66 = 4846000
    67 = $null
    68 = 0
    69 = 0
    70 = 0
    71 = 0
    72 = 2750000
    73 = 0
    74 = 0
    75 = 0
    76 = 3406000
    77 = 0
    80 = 43373
    81 = 114000
    82 = $null
    83 = 312000
    84 = 0
    85 = 0
    86 = 0
    87 = 0
    88 = 0
    89 = 297000
    90 = $null
    91 = -185000
    92 = 0
    93 = 0
    94 = -202000
    95 = $null
    96 = -60000
    97 = 0
    98 = 0
    99 = 0
    100 = -368000
    101 = 1000
    102 = -272000
}

foreach ($row in $colD.Keys) {
    $val = $colD[$row]
    if ($null -ne $val) {
        $ws.Cells.Item($row, 4).Value2 = $val
    }
}

# A handful of rows also got their previously-most-recent period (now in
# column E) corrected/restated at the same time.
$colE = @{
    9 = 2346000
    10 = 914000
    14 = 132000
    17 = 3120000
    18 = 140000
    20 = -6000
    32 = 6000
    89 = 383000
    100 = 119000
}

foreach ($row in $colE.Keys) {
    $ws.Cells.Item($row, 5).Value2 = $colE[$row]
}
